$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 84 (shifts existing rows 84-186 down to 86-188),
# carrying the D-column date formatting down from the row above (Excel default
# behaviour for Insert, matching the inserted <c r="D.." s="2"/> in the diff).
$ws.Range("A84:A85").EntireRow.Insert()

# Row 84: Comercializadora del Agro de Limarí - Poroto verde - Magnum - Primera
$ws.Cells.Item(84, 1).Value = 2
$ws.Cells.Item(84, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(84, 3).Value = "Coquimbo"
$ws.Cells.Item(84, 4).Value = 44664
$ws.Cells.Item(84, 5).Value = 4
$ws.Cells.Item(84, 6).Value = 100112031
$ws.Cells.Item(84, 7).Value = "Poroto verde"
$ws.Cells.Item(84, 8).Value = "Magnum"
$ws.Cells.Item(84, 9).Value = "Primera"
$ws.Cells.Item(84, 10).Value = 400
$ws.Cells.Item(84, 11).Value = 14000
$ws.Cells.Item(84, 12).Value = 15000
$ws.Cells.Item(84, 13).Value = 14500
$ws.Cells.Item(84, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(84, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(84, 16).Value = 580
$ws.Cells.Item(84, 17).Value = 25
$ws.Cells.Item(84, 18).Value = "Hortaliza"

# Row 85: Comercializadora del Agro de Limarí - Poroto verde - Sin especificar - Primera
$ws.Cells.Item(85, 1).Value = 2
$ws.Cells.Item(85, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(85, 3).Value = "Coquimbo"
$ws.Cells.Item(85, 4).Value = 44664
$ws.Cells.Item(85, 5).Value = 4
$ws.Cells.Item(85, 6).Value = 100112031
$ws.Cells.Item(85, 7).Value = "Poroto verde"
$ws.Cells.Item(85, 8).Value = "Sin especificar"
$ws.Cells.Item(85, 9).Value = "Primera"
$ws.Cells.Item(85, 10).Value = 500
$ws.Cells.Item(85, 11).Value = 17000
$ws.Cells.Item(85, 12).Value = 18000
$ws.Cells.Item(85, 13).Value = 17500
$ws.Cells.Item(85, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(85, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(85, 16).Value = 700
$ws.Cells.Item(85, 17).Value = 25
$ws.Cells.Item(85, 18).Value = "Hortaliza"
